# Kings XI Punjab / Murugan Ashwin workbook:
# Append a new row (row 3) to the match-log sheet that duplicates the
# existing row 2 data (same venue/date/result/teams/batsman stats),
# extending the used range from A1:K2 to A1:K3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2's data range into row 3 so the new row carries the exact same
# literal (text) values and cell formatting as row 2 - this is how the
# source data (duplicate match entry) ended up appended to the sheet.
$ws.Range("A2:K2").Copy($ws.Range("A3:K3"))
